# Generate Report for Archive
#
# Status text changes from "Ready for handoff" to "In Translation" on all
# three sheets (Overview summary columns E/F, and the Status column on the
# zh-cn / de-de per-language sheets). Updating the text shortens the
# longest value in those columns, so the sheet's column widths are then
# refreshed to fit the new content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status values.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# Refresh the column widths that depend on the status text length.
$wsOverview.Columns.Item(5).ColumnWidth = 12.576851254417766
$wsOverview.Columns.Item(6).ColumnWidth = 12.576851254417766
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.576851254417766
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.576851254417766
